# "Generate Report for Handoff"
#
# Updates the status of the "b.md" file (row 3) from
# "Handed back: in sync with en-US" to "Ready for handoff" on the
# Overview / zh-cn / de-de sheets, and records the new handoff file name
# + datetime on the two locale detail sheets.

$wb = $excel.ActiveWorkbook

$newHandoffBase = "b.63290e5768f688058c7b37413b0a5c26c308f864"

# Re-adding a hyperlink through the object model re-applies Excel's
# built-in "Hyperlink" look to the target cell. Restore the exact
# font this workbook already uses for hyperlinked cells (Calibri 11,
# single underline, RGB FF6495ED) so the visual style stays the same
# as every other (untouched) hyperlink cell.
function Restore-HyperlinkFont {
    param($Range)
    $f = $Range.Font
    $f.Name = "Calibri"
    $f.Size = 11
    $f.Underline = 2
    $f.Color = 15570276
}

# ---------------------------------------------------------------------
# Sheet 1: "Overview" - just flip the status text for the b.md row
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("C3").Value = "$newHandoffBase.zh-cn.xlf"
$wsZhCn.Range("D3").Value = "2016-03-04 15:46:17"

# Rebuild the hyperlinks collection: every hyperlink keeps the exact same
# target address and display text it had before, except for C3 whose
# display text now shows the new handoff file name (its target address
# is unchanged).
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b63cef3342338c6a3236105c4bfc8a2698d019bc/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cb46c61ea02a08229ea7459dec6fb285e001bc6f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/18c9f07d224370118b32b85e458fcb5afcb213bf/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/8f1cab86e454a977fe7f30a6017db22b6176c1fe/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b63cef3342338c6a3236105c4bfc8a2698d019bc/e2e/b.md", [Type]::Missing, [Type]::Missing, "b.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cb46c61ea02a08229ea7459dec6fb285e001bc6f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "$newHandoffBase.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/18c9f07d224370118b32b85e458fcb5afcb213bf/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/8f1cab86e454a977fe7f30a6017db22b6176c1fe/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/b63cef3342338c6a3236105c4bfc8a2698d019bc/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

foreach ($addr in @("A2","C2","E2","F2","A3","C3","E3","F3","A4")) {
    Restore-HyperlinkFont -Range $wsZhCn.Range($addr)
}

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("C3").Value = "$newHandoffBase.de-de.xlf"
$wsDeDe.Range("D3").Value = "2016-03-04 15:46:30"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/b63cef3342338c6a3236105c4bfc8a2698d019bc/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3c33d6c30a3e9c874c4a4b5caee25fb94dfe8419/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e9498d37640da4cce5445bdae630529ff12beac1/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d12a9ff72ef63fa42838eb0e3a900015d4332991/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/b63cef3342338c6a3236105c4bfc8a2698d019bc/e2e/b.md", [Type]::Missing, [Type]::Missing, "b.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3c33d6c30a3e9c874c4a4b5caee25fb94dfe8419/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", [Type]::Missing, [Type]::Missing, "$newHandoffBase.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e9498d37640da4cce5445bdae630529ff12beac1/e2e/a.md", [Type]::Missing, [Type]::Missing, "a.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d12a9ff72ef63fa42838eb0e3a900015d4332991/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", [Type]::Missing, [Type]::Missing, "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/b63cef3342338c6a3236105c4bfc8a2698d019bc/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

foreach ($addr in @("A2","C2","E2","F2","A3","C3","E3","F3","A4")) {
    Restore-HyperlinkFont -Range $wsDeDe.Range($addr)
}
